# Apply the committed changes to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values in row 4 (E4:K4)
$ws.Range("E4").Value = 1871
$ws.Range("F4").Value = 1871
$ws.Range("G4").Value = 1841
$ws.Range("H4").Value = 1832
$ws.Range("I4").Value = 1964
$ws.Range("J4").Value = 2008
$ws.Range("K4").Value = 2034

# Update the selection: activeCell A3, sqref A3
$ws.Range("A3").Select()
